$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 514.6667
$ws.Range("I6").Value = 448.14285
$ws.Range("K6").Value = 1344.42855
$ws.Range("M6").Value = -1232.42855
$ws.Range("I39").Value = 326
$ws.Range("J39").Value = 1699.6666
$ws.Range("K39").Value = 978
$ws.Range("L39").Value = 5098.9998
$ws.Range("M39").Value = -682
$ws.Range("N39").Value = -5690.9998
$ws.Range("H42").Value = 7.3333335
$ws.Range("I42").Value = 7.3333335
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 22.0000005
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 207.9999995
$ws.Range("N42").ClearContents()
$ws.Range("H86").Value = 5048.3076
$ws.Range("I86").Value = 3828.5715
$ws.Range("J86").Value = 6471.3335
$ws.Range("K86").Value = 3828.5715
$ws.Range("L86").Value = 6471.3335
$ws.Range("M86").Value = -2705.5715
$ws.Range("N86").Value = -8717.333500000001
$ws.Range("H89").Value = 5048.3076
$ws.Range("I89").Value = 3828.5715
$ws.Range("J89").Value = 6471.3335
$ws.Range("K89").Value = 19142.8575
$ws.Range("L89").Value = 32356.6675
$ws.Range("M89").Value = -13526.8575
$ws.Range("N89").Value = -43588.6675
$ws.Range("H137").Value = 875.0714
$ws.Range("I137").Value = 631.5
$ws.Range("J137").Value = 1199.8334
$ws.Range("K137").Value = 1894.5
$ws.Range("L137").Value = 3599.5002
$ws.Range("M137").Value = 655.5
$ws.Range("N137").Value = -8699.5002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 4269
$ws.Range("J30").Value = 4269
$ws.Range("L30").Value = 4269
$ws.Range("N30").Value = -4569
$ws.Range("H32").Value = 2366.3447
$ws.Range("I32").Value = 2173.3704
$ws.Range("K32").Value = 2173.3704
$ws.Range("M32").Value = -1886.3704
$ws.Range("H37").Value = 7799.8
$ws.Range("I37").Value = 7799.8
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 7799.8
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -7526.8
$ws.Range("N37").ClearContents()
$ws.Range("H45").Value = 4119.65
$ws.Range("I45").Value = 1446.5
$ws.Range("K45").Value = 1446.5
$ws.Range("M45").Value = -1069.5
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H102").Value = 1842.1666
$ws.Range("I102").Value = 1910.7
$ws.Range("K102").Value = 1910.7
$ws.Range("M102").Value = -288.7
$ws.Range("H110").Value = 1075
$ws.Range("I110").Value = 1112.5
$ws.Range("K110").Value = 1112.5
$ws.Range("M110").Value = 932.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 5080.6665
$ws.Range("I33").Value = 5080.6665
$ws.Range("K33").Value = 5080.6665
$ws.Range("M33").Value = -4744.6665
$ws.Range("H54").Value = 12100
$ws.Range("I54").Value = 11000
$ws.Range("J54").Value = 16500
$ws.Range("K54").Value = 11000
$ws.Range("L54").Value = 16500
$ws.Range("M54").Value = -10516
$ws.Range("N54").Value = -17468
$ws.Range("H86").Value = 2293.5
$ws.Range("I86").Value = 3933
$ws.Range("J86").Value = 1309.8
$ws.Range("K86").Value = 3933
$ws.Range("L86").Value = 1309.8
$ws.Range("M86").Value = -2810
$ws.Range("N86").Value = -3555.8
$ws.Range("H89").Value = 2293.5
$ws.Range("I89").Value = 3933
$ws.Range("J89").Value = 1309.8
$ws.Range("K89").Value = 19665
$ws.Range("L89").Value = 6549
$ws.Range("M89").Value = -14049
$ws.Range("N89").Value = -17781
$ws.Range("H94").Value = 4680.727
$ws.Range("J94").Value = 4832.5
$ws.Range("L94").Value = 4832.5
$ws.Range("N94").Value = -5734.5
$ws.Range("H99").Value = 4009.5
$ws.Range("I99").Value = 4009.5
$ws.Range("K99").Value = 4009.5
$ws.Range("M99").Value = -2511.5
$ws.Range("H107").Value = 937
$ws.Range("I107").Value = 731.6667
$ws.Range("K107").Value = 731.6667
$ws.Range("M107").Value = 1188.3333
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1999
$ws.Range("I21").Value = 1999
$ws.Range("K21").Value = 1999
$ws.Range("M21").Value = -1764
$ws.Range("H31").Value = 2381
$ws.Range("J31").Value = 2822.25
$ws.Range("L31").Value = 2822.25
$ws.Range("N31").Value = -3412.25
$ws.Range("H34").Value = 2381
$ws.Range("J34").Value = 2822.25
$ws.Range("L34").Value = 2822.25
$ws.Range("N34").Value = -3226.25
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 1870.75
$ws.Range("I122").Value = 1658.8334
$ws.Range("K122").Value = 4976.5002
$ws.Range("M122").Value = -2526.5002
$ws.Range("H134").Value = 2357.6
$ws.Range("I134").Value = 2193.5
$ws.Range("K134").Value = 6580.5
$ws.Range("M134").Value = -4045.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 250000260
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 600
$ws.Range("M4").Value = -488
$ws.Range("H15").Value = 94.5
$ws.Range("J15").Value = 117
$ws.Range("L15").Value = 351
$ws.Range("N15").Value = -631
$ws.Range("H47").Value = 98.5
$ws.Range("I47").Value = 98.5
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 295.5
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 135.5
$ws.Range("N47").ClearContents()
$ws.Range("H69").Value = 3991.3333
$ws.Range("I69").Value = 3990
$ws.Range("K69").Value = 11970
$ws.Range("M69").Value = -11159
$ws.Range("H72").Value = 3991.3333
$ws.Range("I72").Value = 3990
$ws.Range("K72").Value = 35910
$ws.Range("M72").Value = -31854
$ws.Range("H105").Value = 7000
$ws.Range("I105").Value = 7000
$ws.Range("K105").Value = 21000
$ws.Range("M105").Value = -18379

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 44000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 44000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H80").Value = 2971.6
$ws.Range("I80").Value = 851
$ws.Range("K80").Value = 851
$ws.Range("M80").Value = 147
$ws.Range("H83").Value = 2971.6
$ws.Range("I83").Value = 851
$ws.Range("K83").Value = 4255
$ws.Range("M83").Value = 737
$ws.Range("H97").Value = 1799.9
$ws.Range("I97").Value = 1495.5714
$ws.Range("K97").Value = 1495.5714
$ws.Range("M97").Value = -999.5714
$ws.Range("H102").Value = 1912.6471
$ws.Range("I102").Value = 1813.4375
$ws.Range("K102").Value = 1813.4375
$ws.Range("M102").Value = -191.4375
$ws.Range("H126").Value = 5380
$ws.Range("I126").Value = 4966.6665
$ws.Range("K126").Value = 14899.9995
$ws.Range("M126").Value = -12429.9995
$ws.Range("H132").Value = 5666.6665
$ws.Range("I132").Value = 5666.6665
$ws.Range("K132").Value = 16999.9995
$ws.Range("M132").Value = -14469.9995
$ws.Range("H136").Value = 22831.5
$ws.Range("I136").Value = 15000
$ws.Range("J136").Value = 30663
$ws.Range("K136").Value = 45000
$ws.Range("L136").Value = 91989
$ws.Range("M136").Value = -42450
$ws.Range("N136").Value = -97089

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 49388.75
$ws.Range("I63").Value = 46000
$ws.Range("J63").Value = 50518.332
$ws.Range("K63").Value = 46000
$ws.Range("L63").Value = 50518.332
$ws.Range("M63").Value = -45251
$ws.Range("N63").Value = -52016.332
$ws.Range("H66").Value = 49388.75
$ws.Range("I66").Value = 46000
$ws.Range("J66").Value = 50518.332
$ws.Range("K66").Value = 138000
$ws.Range("L66").Value = 151554.996
$ws.Range("M66").Value = -134256
$ws.Range("N66").Value = -159042.996
$ws.Range("H100").Value = 2800
$ws.Range("I100").Value = 2500
$ws.Range("J100").Value = 3400
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 3400
$ws.Range("M100").Value = -1959
$ws.Range("N100").Value = -4482

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5810305
$ws.Range("I100").Value = 8713670
$ws.Range("K100").Value = 17427340
$ws.Range("M100").Value = -17426799
$ws.Range("H112").Value = 30945
$ws.Range("J112").Value = 30945
$ws.Range("L112").Value = 30945
$ws.Range("N112").Value = -33899
